$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F18").Value = 4710
$ws1.Range("F20").Value = 822
$ws1.Range("F21").Value = 105
$ws1.Range("F22").Value = 2201
$ws1.Range("F25").Value = 2067

# Sheet "全部类型" (All Types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F18").Value = 4710
$ws4.Range("F22").Value = 822
$ws4.Range("F23").Value = 105
$ws4.Range("F24").Value = 2201
$ws4.Range("F27").Value = 2067
